# Apply updated numeric values (F/G columns) per commit: 'Update gh-pages to output generated at 456a3b4'
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3253
$ws.Range("F8").Value = 7608
$ws.Range("F11").Value = 9
$ws.Range("F12").Value = 23
$ws.Range("F13").Value = 139
$ws.Range("F14").Value = 650
$ws.Range("F15").Value = 1099
$ws.Range("F16").Value = 1033
$ws.Range("F19").Value = 1405
$ws.Range("G19").Value = 98
$ws.Range("F21").Value = 6020
$ws.Range("F22").Value = 23
$ws.Range("F24").Value = 4161
$ws.Range("F25").Value = 3319
$ws.Range("F27").Value = 86
$ws.Range("F28").Value = 86
$ws.Range("F29").Value = 1024
$ws.Range("F30").Value = 264
$ws.Range("F33").Value = 1022
$ws.Range("F35").Value = 1010
$ws.Range("F36").Value = 77
$ws.Range("F37").Value = 69
$ws.Range("F42").Value = 573
$ws.Range("F43").Value = 365
$ws.Range("F45").Value = 1053
$ws.Range("F46").Value = 462
$ws.Range("F48").Value = 2194
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 595
$ws.Range("F12").Value = 72
$ws.Range("F14").Value = 76
$ws.Range("F15").Value = 174
$ws.Range("F25").Value = 22
$ws.Range("F27").Value = 4584
$ws.Range("F28").Value = 4584
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1939
$ws.Range("F7").Value = 1895
$ws.Range("F10").Value = 1245
$ws.Range("F12").Value = 510
$ws.Range("F13").Value = 2030
$ws.Range("F14").Value = 8707
$ws.Range("F15").Value = 856
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 1939
$ws.Range("F8").Value = 7608
$ws.Range("F10").Value = 1245
$ws.Range("F15").Value = 139
$ws.Range("F16").Value = 856
$ws.Range("F18").Value = 595
$ws.Range("F19").Value = 595
$ws.Range("F20").Value = 650
$ws.Range("F21").Value = 1099
$ws.Range("F22").Value = 1033
$ws.Range("F23").Value = 72
$ws.Range("F26").Value = 174
$ws.Range("F27").Value = 1405
$ws.Range("G27").Value = 98
$ws.Range("F29").Value = 6020
$ws.Range("F31").Value = 4161
$ws.Range("F32").Value = 3319
$ws.Range("F33").Value = 86
$ws.Range("F34").Value = 1024
$ws.Range("F35").Value = 264
$ws.Range("F36").Value = 1022
$ws.Range("F38").Value = 77
$ws.Range("F39").Value = 69
$ws.Range("F42").Value = 573
$ws.Range("F43").Value = 365
$ws.Range("F45").Value = 22
$ws.Range("F46").Value = 462
$ws.Range("F47").Value = 2194
$ws.Range("F49").Value = 4584
